$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.715.40"
$ws.Range("E2").Value = "  +3.05%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.447.92"
$ws.Range("E3").Value = "  +1.94%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.83%  "

# Row 6 - Solana
$ws.Range("E6").Value = "  +2.89%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.56%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.447.01"
$ws.Range("E9").Value = "  +1.60%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.53%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.49%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +1.31%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +2.98%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.49%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +5.76%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("E16").Value = "  +3.79%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "62.605.95"
$ws.Range("E17").Value = "  +3.10%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.449.93"
$ws.Range("E18").Value = "  +1.60%  "

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.09%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +2.78%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.40%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  +1.14%  "

# Row 23 - SuiNetwork
$ws.Range("E23").Value = "  +10.52%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.16%  "

# Row 25 - Litecoin
$ws.Range("E25").Value = "  +1.16%  "

# Row 26 - Bittensor
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "643.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.69%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.93%  "

# Row 28 - Aptos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.10%  "

# Row 29 - row29 (BabyDogeCoin -> PEPE)
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0987"
$ws.Range("E29").Value = "  +5.08%  "

# Row 30 - row30 (PEPE -> BabyDogeCoin)
$ws.Range("B30").Value = "BabyDogeCoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D30").Value = "0.0₆0514"
$ws.Range("E30").Value = "  +84.57%  "

# Row 31 - WrappedeETH
$ws.Range("D31").Value = "2.575.17"
$ws.Range("E31").Value = "  +2.28%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  +8.29%  "

# Row 34 - PancakeSwap
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.59%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +4.98%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +2.27%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.18%  "

# Row 38 - NEARProtocol
$ws.Range("E38").Value = "  +3.55%  "

# Row 39 - RenderToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.17%  "

# Row 40 - Monero
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "153.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "

# Row 41 - PolygonEcosystemToken
$ws.Range("E41").Value = "  +0.99%  "

# Row 42 - EthereumClassic
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.91%  "

# Row 43 - dogwifhat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.37%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  +5.31%  "

# Row 45 - row45 (USDe -> OKB)
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.66%  "

# Row 46 - row46 (OKB -> USDe)
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "

# Row 47 - WhiteBITCoin
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "14.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +27.72%  "

# Row 48 - Aave
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "145.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.25%  "

# Row 49 - Filecoin
$ws.Range("E49").Value = "  +2.54%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.06%  "

# Row 51 - Mantle
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.605"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.96%  "

